# In the "Recorded By" column (G), swap the order of "System" and the
# recorder's email address wherever both appear together, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#
# Cells that contain only "System", only the email address, or any other
# text are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G ("Recorded By")
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
